$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.00469825789230973
$ws.Range("J2").Value = 0.00469825789230973
$ws.Range("M2").Value = 45.90594266666667
$ws.Range("N2").Value = 137.717828
$ws.Range("O2").Value = 0.3954672001633582
$ws.Range("P2").Value = 0.3954672001633583
$ws.Range("Q2").Value = 4.288915713442223
$ws.Range("R2").Value = 38.60024142098
$ws.Range("S2").Value = 0.001858006894317129
$ws.Range("T2").Value = 0.00185800689431713
$ws.Range("I3").Value = 0.00469825789230973
$ws.Range("J3").Value = 0.00469825789230973
$ws.Range("O3").Value = 0.3484294080560655
$ws.Range("P3").Value = 0.3484294080560656
$ws.Range("S3").Value = 0.001637011216312218
$ws.Range("T3").Value = 0.001637011216312218
$ws.Range("I4").Value = 0.00469825789230973
$ws.Range("J4").Value = 0.00469825789230973
$ws.Range("M4").Value = 12.761795
$ws.Range("N4").Value = 38.28538500000001
$ws.Range("O4").Value = 0.1099393900775594
$ws.Range("P4").Value = 0.1099393900775594
$ws.Range("Q4").Value = 1.192313237191667
$ws.Range("R4").Value = 10.730819134725
$ws.Range("S4").Value = 0.0005165236071076115
$ws.Range("T4").Value = 0.0005165236071076116
$ws.Range("I5").Value = 0.00469825789230973
$ws.Range("J5").Value = 0.00469825789230973
$ws.Range("M5").Value = 16.966758
$ws.Range("N5").Value = 50.900274
$ws.Range("O5").Value = 0.1461640017030168
$ws.Range("P5").Value = 0.1461640017030168
$ws.Range("Q5").Value = 1.58517592201
$ws.Range("R5").Value = 14.26658329809
$ws.Range("S5").Value = 0.0006867161745727716
$ws.Range("T5").Value = 0.0006867161745727717
$ws.Range("I6").Value = 0.7185612021237531
$ws.Range("J6").Value = 0.7185612021237531
$ws.Range("M6").Value = 45.90594266666667
$ws.Range("N6").Value = 137.717828
$ws.Range("O6").Value = 0.3954672001633582
$ws.Range("P6").Value = 0.3954672001633583
$ws.Range("Q6").Value = 655.9555693830627
$ws.Range("R6").Value = 5903.600124447564
$ws.Range("S6").Value = 0.2841673867498976
$ws.Range("T6").Value = 0.2841673867498976
$ws.Range("I7").Value = 0.7185612021237531
$ws.Range("J7").Value = 0.7185612021237531
$ws.Range("O7").Value = 0.3484294080560655
$ws.Range("P7").Value = 0.3484294080560656
$ws.Range("Q7").Value = 577.9346824637026
$ws.Range("R7").Value = 5201.412142173323
$ws.Range("S7").Value = 0.2503678543080342
$ws.Range("T7").Value = 0.2503678543080342
$ws.Range("I8").Value = 0.7185612021237531
$ws.Range("J8").Value = 0.7185612021237531
$ws.Range("M8").Value = 12.761795
$ws.Range("N8").Value = 38.28538500000001
$ws.Range("O8").Value = 0.1099393900775594
$ws.Range("P8").Value = 0.1099393900775594
$ws.Range("Q8").Value = 182.354832932195
$ws.Range("R8").Value = 1641.193496389755
$ws.Range("S8").Value = 0.0789981802948833
$ws.Range("T8").Value = 0.07899818029488331
$ws.Range("I9").Value = 0.7185612021237531
$ws.Range("J9").Value = 0.7185612021237531
$ws.Range("M9").Value = 16.966758
$ws.Range("N9").Value = 50.900274
$ws.Range("O9").Value = 0.1461640017030168
$ws.Range("P9").Value = 0.1461640017030168
$ws.Range("Q9").Value = 242.440058039718
$ws.Range("R9").Value = 2181.960522357462
$ws.Range("S9").Value = 0.105027780770938
$ws.Range("T9").Value = 0.1050277807709381
$ws.Range("G10").Value = 5.503190333333333
$ws.Range("H10").Value = 16.509571
$ws.Range("I10").Value = 0.2767405399839373
$ws.Range("J10").Value = 0.2767405399839373
$ws.Range("M10").Value = 45.90594266666667
$ws.Range("N10").Value = 137.717828
$ws.Range("O10").Value = 0.3954672001633582
$ws.Range("P10").Value = 0.3954672001633583
$ws.Range("Q10").Value = 252.6291399257542
$ws.Range("R10").Value = 2273.662259331788
$ws.Range("S10").Value = 0.1094418065191436
$ws.Range("T10").Value = 0.1094418065191436
$ws.Range("G11").Value = 5.503190333333333
$ws.Range("H11").Value = 16.509571
$ws.Range("I11").Value = 0.2767405399839373
$ws.Range("J11").Value = 0.2767405399839373
$ws.Range("O11").Value = 0.3484294080560655
$ws.Range("P11").Value = 0.3484294080560656
$ws.Range("Q11").Value = 222.5808401953009
$ws.Range("R11").Value = 2003.227561757708
$ws.Range("S11").Value = 0.0964245425317192
$ws.Range("T11").Value = 0.09642454253171921
$ws.Range("G12").Value = 5.503190333333333
$ws.Range("H12").Value = 16.509571
$ws.Range("I12").Value = 0.2767405399839373
$ws.Range("J12").Value = 0.2767405399839373
$ws.Range("M12").Value = 12.761795
$ws.Range("N12").Value = 38.28538500000001
$ws.Range("O12").Value = 0.1099393900775594
$ws.Range("P12").Value = 0.1099393900775594
$ws.Range("Q12").Value = 70.23058687998167
$ws.Range("R12").Value = 632.0752819198351
$ws.Range("S12").Value = 0.0304246861755685
$ws.Range("T12").Value = 0.03042468617556851
$ws.Range("G13").Value = 5.503190333333333
$ws.Range("H13").Value = 16.509571
$ws.Range("I13").Value = 0.2767405399839373
$ws.Range("J13").Value = 0.2767405399839373
$ws.Range("M13").Value = 16.966758
$ws.Range("N13").Value = 50.900274
$ws.Range("O13").Value = 0.1461640017030168
$ws.Range("P13").Value = 0.1461640017030168
$ws.Range("Q13").Value = 93.371298613606
$ws.Range("R13").Value = 840.341687522454
$ws.Range("S13").Value = 0.040449504757506
$ws.Range("T13").Value = 0.04044950475750601

Write-Output "Updated cells with new TPM values"
